$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.184.73"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "2.247.51"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.23"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -3.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.14"
$ws.Range("E7").Value = "  +3.82%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.14"
$ws.Range("E10").Value = "  +5.95%  "
$ws.Range("E12").Value = "  -5.00%  "
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "2.583.81"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.76"
$ws.Range("E15").Value = "  -4.29%  "
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "2.257.62"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "42.063.47"
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.94"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.10"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.27"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.91"
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.32"
$ws.Range("E25").Value = "  -2.94%  "
$ws.Range("E26").Value = "  -7.64%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.30"
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  +18.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  +4.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.90"
$ws.Range("E30").Value = "  +1.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.57"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.60"
$ws.Range("E33").Value = "  +5.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.121"
$ws.Range("E34").Value = "  -5.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.51"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.95"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14.49"
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("E41").Value = "  -7.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "112.67"
$ws.Range("E42").Value = "  +7.98%  "
$ws.Range("E43").Value = "  -7.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.89"
$ws.Range("E44").Value = "  -2.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.69"
$ws.Range("E45").Value = "  -5.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0991"
$ws.Range("E46").Value = "  -4.69%  "
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -12.22%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.438"
$ws.Range("E51").Value = "  +12.83%  "
